# Auto-generated Excel COM-interop script applying numeric updates
# to the "Leve Profits" sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# for columns H-N (market price / profit calculations), per scheduled
# market-data runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 8928997
$ws.Range("I19").Value = 28571562
$ws.Range("J19").Value = 559.2727
$ws.Range("K19").Value = 28571562
$ws.Range("L19").Value = 559.2727
$ws.Range("M19").Value = -28571387
$ws.Range("N19").Value = -909.2727
$ws.Range("H62").Value = 2641.7222
$ws.Range("I62").Value = 2122.2222
$ws.Range("J62").Value = 3161.2222
$ws.Range("K62").Value = 2122.2222
$ws.Range("L62").Value = 3161.2222
$ws.Range("M62").Value = -1498.2222
$ws.Range("N62").Value = -4409.2222
$ws.Range("H65").Value = 2641.7222
$ws.Range("I65").Value = 2122.2222
$ws.Range("J65").Value = 3161.2222
$ws.Range("K65").Value = 10611.111
$ws.Range("L65").Value = 15806.111
$ws.Range("M65").Value = -7491.111000000001
$ws.Range("N65").Value = -22046.111
$ws.Range("H101").Value = 4834
$ws.Range("J101").Value = 9333.333000000001
$ws.Range("L101").Value = 27999.999
$ws.Range("N101").Value = -31243.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 150004000
$ws.Range("I2").Value = 250000670
$ws.Range("K2").Value = 250000670
$ws.Range("M2").Value = -250000557
$ws.Range("H32").Value = 15828.609
$ws.Range("I32").Value = 13759.4375
$ws.Range("J32").Value = 22036.125
$ws.Range("K32").Value = 13759.4375
$ws.Range("L32").Value = 22036.125
$ws.Range("M32").Value = -13472.4375
$ws.Range("N32").Value = -22610.125
$ws.Range("H45").Value = 1428.0714
$ws.Range("I45").Value = 1114.8462
$ws.Range("J45").Value = 5500
$ws.Range("K45").Value = 1114.8462
$ws.Range("L45").Value = 5500
$ws.Range("M45").Value = -737.8462
$ws.Range("N45").Value = -6254
$ws.Range("H63").Value = 2575
$ws.Range("I63").Value = 2450
$ws.Range("J63").Value = 2700
$ws.Range("K63").Value = 2450
$ws.Range("L63").Value = 2700
$ws.Range("M63").Value = -1764
$ws.Range("N63").Value = -4072
$ws.Range("H66").Value = 2575
$ws.Range("I66").Value = 2450
$ws.Range("J66").Value = 2700
$ws.Range("K66").Value = 12250
$ws.Range("L66").Value = 13500
$ws.Range("M66").Value = -8818
$ws.Range("N66").Value = -20364
$ws.Range("H110").Value = 1682.56
$ws.Range("I110").Value = 624
$ws.Range("J110").Value = 3932
$ws.Range("K110").Value = 624
$ws.Range("L110").Value = 3932
$ws.Range("M110").Value = 1421
$ws.Range("N110").Value = -8022
$ws.Range("H116").Value = 150004000
$ws.Range("I116").Value = 250000670
$ws.Range("K116").Value = 250000670
$ws.Range("M116").Value = -249998376
$ws.Range("H132").Value = 2041.7377
$ws.Range("I132").Value = 1492.4681
$ws.Range("K132").Value = 4477.4043
$ws.Range("M132").Value = -1947.4043

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 150004000
$ws.Range("I3").Value = 250000670
$ws.Range("K3").Value = 250000670
$ws.Range("M3").Value = -250000556
$ws.Range("H80").Value = 966.9375
$ws.Range("I80").Value = 1252.6
$ws.Range("J80").Value = 837.0909
$ws.Range("K80").Value = 1252.6
$ws.Range("L80").Value = 837.0909
$ws.Range("M80").Value = -254.5999999999999
$ws.Range("N80").Value = -2833.0909
$ws.Range("H83").Value = 966.9375
$ws.Range("I83").Value = 1252.6
$ws.Range("J83").Value = 837.0909
$ws.Range("K83").Value = 6263
$ws.Range("L83").Value = 4185.4545
$ws.Range("M83").Value = -1271
$ws.Range("N83").Value = -14169.4545
$ws.Range("H86").Value = 28497.37
$ws.Range("I86").Value = 1307.1428
$ws.Range("K86").Value = 1307.1428
$ws.Range("M86").Value = -184.1428000000001
$ws.Range("H89").Value = 28497.37
$ws.Range("I89").Value = 1307.1428
$ws.Range("K89").Value = 6535.714
$ws.Range("M89").Value = -919.7139999999999
$ws.Range("H105").Value = 1683.4445
$ws.Range("I105").Value = 1323.0526
$ws.Range("J105").Value = 2539.375
$ws.Range("K105").Value = 1323.0526
$ws.Range("L105").Value = 2539.375
$ws.Range("M105").Value = 423.9474
$ws.Range("N105").Value = -6033.375
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 11005.2
$ws.Range("J25").Value = 11005.2
$ws.Range("L25").Value = 11005.2
$ws.Range("N25").Value = -11353.2
$ws.Range("H31").Value = 2965.774
$ws.Range("I31").Value = 2004.0638
$ws.Range("J31").Value = 4187.4053
$ws.Range("K31").Value = 2004.0638
$ws.Range("L31").Value = 4187.4053
$ws.Range("M31").Value = -1709.0638
$ws.Range("N31").Value = -4777.4053
$ws.Range("H34").Value = 2965.774
$ws.Range("I34").Value = 2004.0638
$ws.Range("J34").Value = 4187.4053
$ws.Range("K34").Value = 2004.0638
$ws.Range("L34").Value = 4187.4053
$ws.Range("M34").Value = -1802.0638
$ws.Range("N34").Value = -4591.4053
$ws.Range("H122").Value = 2426.0881
$ws.Range("I122").Value = 2109.7778
$ws.Range("J122").Value = 3646.1428
$ws.Range("K122").Value = 6329.3334
$ws.Range("L122").Value = 10938.4284
$ws.Range("M122").Value = -3879.3334
$ws.Range("N122").Value = -15838.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 444.0625
$ws.Range("J12").Value = 695.4
$ws.Range("L12").Value = 2086.2
$ws.Range("N12").Value = -2432.2
$ws.Range("H98").Value = 239.05263
$ws.Range("J98").Value = 323.16666
$ws.Range("L98").Value = 969.4999799999999
$ws.Range("N98").Value = -3965.49998
$ws.Range("H131").Value = 1603.8776
$ws.Range("J131").Value = 1252.1428
$ws.Range("L131").Value = 3756.4284
$ws.Range("N131").Value = -13836.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 42562.72
$ws.Range("I102").Value = 1953.25
$ws.Range("K102").Value = 1953.25
$ws.Range("M102").Value = -331.25
$ws.Range("H126").Value = 773409.6
$ws.Range("I126").Value = 3102.2
$ws.Range("J126").Value = 1254851.8
$ws.Range("K126").Value = 9306.599999999999
$ws.Range("L126").Value = 3764555.4
$ws.Range("M126").Value = -6836.599999999999
$ws.Range("N126").Value = -3769495.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 5732.5
$ws.Range("J45").Value = 5732.5
$ws.Range("L45").Value = 5732.5
$ws.Range("N45").Value = -6714.5
$ws.Range("H108").Value = 39146.668
$ws.Range("J108").Value = 39146.668
$ws.Range("L108").Value = 39146.668
$ws.Range("N108").Value = -46826.668
$ws.Range("H126").Value = 42060
$ws.Range("I126").Value = 61392.53
$ws.Range("K126").Value = 184177.59
$ws.Range("M126").Value = -181707.59
$ws.Range("H132").Value = 18446.934
$ws.Range("I132").Value = 5278.478
$ws.Range("J132").Value = 61714.715
$ws.Range("K132").Value = 15835.434
$ws.Range("L132").Value = 185144.145
$ws.Range("M132").Value = -13305.434
$ws.Range("N132").Value = -190204.145

Write-Host "Applied leve-profit updates across 8 sheets"
